# TimeSheet.xlsx final edits: fill in Feb (rows 21-33) and Apr (rows 14-35)
# activity/hours data, and update sheet-view selections / active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Feb sheet: rows 21-33, columns C..H (3x [Activity Code, Hours] pairs)
# ---------------------------------------------------------------------------
$wsFeb = $wb.Worksheets.Item("Feb")

$febData = @(
    @{R=21; C="C"; V="P06"},
    @{R=21; C="D"; V="OFF"},
    @{R=21; C="E"; V="P06"},
    @{R=21; C="F"; V="OFF"},
    @{R=21; C="G"; V="P06"},
    @{R=21; C="H"; V="OFF"},
    @{R=22; C="C"; V="P06"},
    @{R=22; C="D"; V=1},
    @{R=22; C="E"; V="P06"},
    @{R=22; C="F"; V="OFF"},
    @{R=22; C="G"; V="P06"},
    @{R=22; C="H"; V=1},
    @{R=23; C="C"; V="P06"},
    @{R=23; C="D"; V="OFF"},
    @{R=23; C="E"; V="P06"},
    @{R=23; C="F"; V=2},
    @{R=23; C="G"; V="P06"},
    @{R=23; C="H"; V="OFF"},
    @{R=24; C="C"; V="P06"},
    @{R=24; C="D"; V="OFF"},
    @{R=24; C="E"; V="P06"},
    @{R=24; C="F"; V=1},
    @{R=24; C="G"; V="P06"},
    @{R=24; C="H"; V=2},
    @{R=25; C="C"; V="P06"},
    @{R=25; C="D"; V=2},
    @{R=25; C="E"; V="P06"},
    @{R=25; C="F"; V=1.5},
    @{R=25; C="G"; V="P06"},
    @{R=25; C="H"; V=1},
    @{R=26; C="C"; V="P06"},
    @{R=26; C="D"; V="OFF"},
    @{R=26; C="E"; V="P06"},
    @{R=26; C="F"; V="OFF"},
    @{R=26; C="G"; V="P06"},
    @{R=26; C="H"; V="OFF"},
    @{R=27; C="C"; V="P06"},
    @{R=27; C="D"; V=1},
    @{R=27; C="E"; V="P06"},
    @{R=27; C="F"; V="OFF"},
    @{R=27; C="G"; V="P06"},
    @{R=27; C="H"; V=0.5},
    @{R=28; C="C"; V="P06"},
    @{R=28; C="D"; V="OFF"},
    @{R=28; C="E"; V="P06"},
    @{R=28; C="F"; V=1},
    @{R=28; C="G"; V="P06"},
    @{R=28; C="H"; V="OFF"},
    @{R=29; C="C"; V="P06"},
    @{R=29; C="D"; V=1},
    @{R=29; C="E"; V="P06"},
    @{R=29; C="F"; V="OFF"},
    @{R=29; C="G"; V="P06"},
    @{R=29; C="H"; V=1},
    @{R=30; C="C"; V="P07"},
    @{R=30; C="D"; V=1},
    @{R=30; C="E"; V="P07"},
    @{R=30; C="F"; V="OFF"},
    @{R=30; C="G"; V="P07"},
    @{R=30; C="H"; V=1},
    @{R=31; C="C"; V="P07"},
    @{R=31; C="D"; V=1},
    @{R=31; C="E"; V="P07"},
    @{R=31; C="F"; V=1},
    @{R=31; C="G"; V="P07"},
    @{R=31; C="H"; V=1},
    @{R=32; C="C"; V="P07"},
    @{R=32; C="D"; V="OFF"},
    @{R=32; C="E"; V="P07"},
    @{R=32; C="F"; V=2},
    @{R=32; C="G"; V="P07"},
    @{R=32; C="H"; V=1},
    @{R=33; C="C"; V="P07"},
    @{R=33; C="D"; V=1},
    @{R=33; C="E"; V="P07"},
    @{R=33; C="F"; V="OFF"},
    @{R=33; C="G"; V="P07"},
    @{R=33; C="H"; V=1}
)

foreach ($cell in $febData) {
    $wsFeb.Range("$($cell.C)$($cell.R)").Value = $cell.V
}

# ---------------------------------------------------------------------------
# Apr sheet: rows 14-35, columns C..H
# ---------------------------------------------------------------------------
$wsApr = $wb.Worksheets.Item("Apr")

# Rows 14-22 pick up column E's fill from column C (style shifts 12 -> 13)
# before the values are written, matching the source workbook exactly.
for ($r = 14; $r -le 22; $r++) {
    $wsApr.Cells.Item($r, 3).Copy($wsApr.Cells.Item($r, 5))
}

$aprData = @(
    @{R=14; C="C"; V="P09"},
    @{R=14; C="D"; V=1},
    @{R=14; C="E"; V="P09"},
    @{R=14; C="F"; V=1},
    @{R=14; C="G"; V="P09"},
    @{R=14; C="H"; V=1},
    @{R=15; C="C"; V="P09"},
    @{R=15; C="D"; V="OFF"},
    @{R=15; C="E"; V="P09"},
    @{R=15; C="F"; V=1},
    @{R=15; C="G"; V="P09"},
    @{R=15; C="H"; V="OFF"},
    @{R=16; C="C"; V="P09"},
    @{R=16; C="D"; V=1},
    @{R=16; C="E"; V="P09"},
    @{R=16; C="F"; V="OFF"},
    @{R=16; C="G"; V="P09"},
    @{R=16; C="H"; V=1},
    @{R=17; C="C"; V="P09"},
    @{R=17; C="D"; V="OFF"},
    @{R=17; C="E"; V="P09"},
    @{R=17; C="F"; V="OFF"},
    @{R=17; C="G"; V="P09"},
    @{R=17; C="H"; V=1},
    @{R=18; C="C"; V="P09"},
    @{R=18; C="D"; V=2},
    @{R=18; C="E"; V="P09"},
    @{R=18; C="F"; V=2},
    @{R=18; C="G"; V="P09"},
    @{R=18; C="H"; V="OFF"},
    @{R=19; C="C"; V="P09"},
    @{R=19; C="D"; V="OFF"},
    @{R=19; C="E"; V="P09"},
    @{R=19; C="F"; V="OFF"},
    @{R=19; C="G"; V="P09"},
    @{R=19; C="H"; V=1},
    @{R=20; C="C"; V="P09"},
    @{R=20; C="D"; V=1},
    @{R=20; C="E"; V="P09"},
    @{R=20; C="F"; V="OFF"},
    @{R=20; C="G"; V="P09"},
    @{R=20; C="H"; V="OFF"},
    @{R=21; C="C"; V="P09"},
    @{R=21; C="D"; V="OFF"},
    @{R=21; C="E"; V="P09"},
    @{R=21; C="F"; V=1},
    @{R=21; C="G"; V="P09"},
    @{R=21; C="H"; V=1},
    @{R=22; C="C"; V="P09"},
    @{R=22; C="D"; V=1},
    @{R=22; C="E"; V="P09"},
    @{R=22; C="F"; V="OFF"},
    @{R=22; C="G"; V="P09"},
    @{R=22; C="H"; V="OFF"},
    @{R=23; C="C"; V="P10"},
    @{R=23; C="D"; V=1},
    @{R=23; C="E"; V="P10"},
    @{R=23; C="F"; V=1},
    @{R=23; C="G"; V="P10"},
    @{R=23; C="H"; V=1},
    @{R=24; C="C"; V="P10"},
    @{R=24; C="D"; V="OFF"},
    @{R=24; C="E"; V="P10"},
    @{R=24; C="F"; V="OFF"},
    @{R=24; C="G"; V="P10"},
    @{R=24; C="H"; V="OFF"},
    @{R=25; C="C"; V="P10"},
    @{R=25; C="D"; V="OFF"},
    @{R=25; C="E"; V="P10"},
    @{R=25; C="F"; V=1},
    @{R=25; C="G"; V="P10"},
    @{R=25; C="H"; V="OFF"},
    @{R=26; C="C"; V="P11"},
    @{R=26; C="D"; V=1},
    @{R=26; C="E"; V="P11"},
    @{R=26; C="F"; V="OFF"},
    @{R=26; C="G"; V="P11"},
    @{R=26; C="H"; V="OFF"},
    @{R=27; C="C"; V="P11"},
    @{R=27; C="D"; V="OFF"},
    @{R=27; C="E"; V="P11"},
    @{R=27; C="F"; V=2},
    @{R=27; C="G"; V="P11"},
    @{R=27; C="H"; V=1},
    @{R=28; C="C"; V="P11"},
    @{R=28; C="D"; V=2},
    @{R=28; C="E"; V="P11"},
    @{R=28; C="F"; V="OFF"},
    @{R=28; C="G"; V="P11"},
    @{R=28; C="H"; V="OFF"},
    @{R=29; C="C"; V="P11"},
    @{R=29; C="D"; V="OFF"},
    @{R=29; C="E"; V="P11"},
    @{R=29; C="F"; V="OFF"},
    @{R=29; C="G"; V="P11"},
    @{R=29; C="H"; V=1},
    @{R=30; C="C"; V="P11"},
    @{R=30; C="D"; V="OFF"},
    @{R=30; C="E"; V="P11"},
    @{R=30; C="F"; V=1},
    @{R=30; C="G"; V="P11"},
    @{R=30; C="H"; V="OFF"},
    @{R=31; C="C"; V="P11"},
    @{R=31; C="D"; V="OFF"},
    @{R=31; C="E"; V="P11"},
    @{R=31; C="F"; V="OFF"},
    @{R=31; C="G"; V="P11"},
    @{R=31; C="H"; V=1},
    @{R=32; C="C"; V="NA"},
    @{R=32; C="E"; V="NA"},
    @{R=32; C="G"; V="NA"},
    @{R=33; C="C"; V="NA"},
    @{R=33; C="E"; V="NA"},
    @{R=33; C="G"; V="NA"},
    @{R=34; C="C"; V="NA"},
    @{R=34; C="E"; V="NA"},
    @{R=34; C="G"; V="NA"},
    @{R=35; C="C"; V="NA"},
    @{R=35; C="E"; V="NA"},
    @{R=35; C="G"; V="NA"}
)

foreach ($cell in $aprData) {
    $wsApr.Range("$($cell.C)$($cell.R)").Value = $cell.V
}

# ---------------------------------------------------------------------------
# Sheet-view selections (applied before the final sheet activation below)
# ---------------------------------------------------------------------------
$wsApr.Range("C33:C35").Select()

$wsMar = $wb.Worksheets.Item("Mar")
$wsMar.Range("H35").Select()

# Feb becomes the active/selected tab, with D33 as the active cell - this
# also updates workbook.xml's activeTab and clears tabSelected on Mar/Apr.
$wsFeb.Activate()
$wsFeb.Range("D33").Select()

Write-Output "TimeSheet updated"
